$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (shows up in Overview!B2:C3 and in the zh-cn / de-de "Status" column)
# ---------------------------------------------------------------------------
foreach ($sheet in $wb.Worksheets) {
    $sheet.Cells.Replace("Ready for handoff", "Handed back: in sync with en-US")
}

# Hyperlink font used throughout the workbook for linked filenames
# (underline + RGB 6495ED, expressed as BGR for the COM Color property)
$HyperlinkUnderline = -4119   # xlUnderlineStyleSingle
$HyperlinkColor = 15570276    # BGR for RGB(100,149,237) == FF6495ED

function Set-HyperlinkLook($rng) {
    $rng.Font.Name = "Calibri"
    $rng.Font.Color = $HyperlinkColor
    $rng.Font.Underline = $HyperlinkUnderline
}

# ---------------------------------------------------------------------------
# 2. zh-cn sheet: populate "Latest Target File" (F) / "Latest Handback File"
#    (G) for both data rows, and fill in the real handback timestamp (H).
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Hyperlinks.Add(
    $zh.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/0f29589865e7073b5d7e52f89debd4beae6d441d/e2e/0e9f89cf-6329-44b1-bee3-047f72061dab.md",
    "",
    "",
    "0e9f89cf-6329-44b1-bee3-047f72061dab.md"
)
Set-HyperlinkLook $zh.Range("F2")

$zh.Hyperlinks.Add(
    $zh.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8a45c9af54f747648f9dbc73d95fad9c7f42f940/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/0e9f89cf-6329-44b1-bee3-047f72061dab.a2162950077b55203a9c95d943c441e3793fc4c8.zh-cn.xlf",
    "",
    "",
    "0e9f89cf-6329-44b1-bee3-047f72061dab.a2162950077b55203a9c95d943c441e3793fc4c8.zh-cn.xlf"
)
Set-HyperlinkLook $zh.Range("G2")

$zh.Hyperlinks.Add(
    $zh.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/0f29589865e7073b5d7e52f89debd4beae6d441d/e2e/dd20ad19-4979-4479-97e5-e9c3df9306fe.md",
    "",
    "",
    "dd20ad19-4979-4479-97e5-e9c3df9306fe.md"
)
Set-HyperlinkLook $zh.Range("F3")

$zh.Hyperlinks.Add(
    $zh.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8a45c9af54f747648f9dbc73d95fad9c7f42f940/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/dd20ad19-4979-4479-97e5-e9c3df9306fe.3a38cc87ad4cfedab8990032b3c0d8c5672cc2b5.zh-cn.xlf",
    "",
    "",
    "dd20ad19-4979-4479-97e5-e9c3df9306fe.3a38cc87ad4cfedab8990032b3c0d8c5672cc2b5.zh-cn.xlf"
)
Set-HyperlinkLook $zh.Range("G3")

$zh.Range("H2").Value = "2016-03-19 04:39:52"
$zh.Range("H3").Value = "2016-03-19 04:39:52"

# ---------------------------------------------------------------------------
# 3. de-de sheet: same shape of change, different timestamp/hashes.
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Hyperlinks.Add(
    $de.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/0f29589865e7073b5d7e52f89debd4beae6d441d/e2e/0e9f89cf-6329-44b1-bee3-047f72061dab.md",
    "",
    "",
    "0e9f89cf-6329-44b1-bee3-047f72061dab.md"
)
Set-HyperlinkLook $de.Range("F2")

$de.Hyperlinks.Add(
    $de.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c894cbf859c46b8c43c48c8d385dd1839c0289a2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/0e9f89cf-6329-44b1-bee3-047f72061dab.a2162950077b55203a9c95d943c441e3793fc4c8.de-de.xlf",
    "",
    "",
    "0e9f89cf-6329-44b1-bee3-047f72061dab.a2162950077b55203a9c95d943c441e3793fc4c8.de-de.xlf"
)
Set-HyperlinkLook $de.Range("G2")

$de.Hyperlinks.Add(
    $de.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/0f29589865e7073b5d7e52f89debd4beae6d441d/e2e/dd20ad19-4979-4479-97e5-e9c3df9306fe.md",
    "",
    "",
    "dd20ad19-4979-4479-97e5-e9c3df9306fe.md"
)
Set-HyperlinkLook $de.Range("F3")

$de.Hyperlinks.Add(
    $de.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c894cbf859c46b8c43c48c8d385dd1839c0289a2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/dd20ad19-4979-4479-97e5-e9c3df9306fe.3a38cc87ad4cfedab8990032b3c0d8c5672cc2b5.de-de.xlf",
    "",
    "",
    "dd20ad19-4979-4479-97e5-e9c3df9306fe.3a38cc87ad4cfedab8990032b3c0d8c5672cc2b5.de-de.xlf"
)
Set-HyperlinkLook $de.Range("G3")

$de.Range("H2").Value = "2016-03-19 04:39:58"
$de.Range("H3").Value = "2016-03-19 04:39:58"
